# Update F-column ("想去人数" / want-to-go count) values on sheets
# "展览" (Exhibition), "演出" (Performance) and "全部类型" (All types)
# to match the newly scraped data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 94
$ws1.Range("F4").Value = 38
$ws1.Range("F5").Value = 0
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 6656
$ws1.Range("F10").Value = 0
$ws1.Range("F11").Value = 15
$ws1.Range("F12").Value = 107
$ws1.Range("F14").Value = 138
$ws1.Range("F15").Value = 17
$ws1.Range("F19").Value = 0
$ws1.Range("F21").Value = 120
$ws1.Range("F22").Value = 452
$ws1.Range("F23").Value = 212
$ws1.Range("F24").Value = 195

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 0

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 94
$ws4.Range("F4").Value = 38
$ws4.Range("F5").Value = 450
$ws4.Range("F7").Value = 6656
$ws4.Range("F10").Value = 1290
$ws4.Range("F12").Value = 107
$ws4.Range("F13").Value = 0
$ws4.Range("F15").Value = 17
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 9
$ws4.Range("F20").Value = 0
$ws4.Range("F22").Value = 99
$ws4.Range("F24").Value = 452
$ws4.Range("F26").Value = 195

Write-Host "Applied F-column updates to sheets 1, 2 and 4."
